$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "27.583.14"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "1.666.73"
$ws.Range("E3").Value = "  -3.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue "D8" "23.57"
$ws.Range("E8").Value = "  -2.23%  "
Set-TextValue "D9" "0.264"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -1.84%  "
Set-TextValue "D11" "0.0882"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "1.902.54"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("D13").Value = "1.661.24"
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("E15").Value = "  -0.77%  "
Set-TextValue "D16" "66.10"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.602.22"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D18" "242.09"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "0.0₃0729"
$ws.Range("E19").Value = "  -3.50%  "
Set-TextValue "D20" "7.52"
$ws.Range("E20").Value = "  -4.40%  "
$ws.Range("E21").Value = "  +0.03%  "
Set-TextValue "D22" "4.48"
$ws.Range("E22").Value = "  -3.16%  "
Set-TextValue "D23" "9.29"
$ws.Range("E23").Value = "  -4.71%  "
$ws.Range("E24").Value = "  -3.67%  "
Set-TextValue "D25" "146.06"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("E26").Value = "  -4.51%  "
Set-TextValue "D27" "16.35"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("E30").Value = "  +3.89%  "
Set-TextValue "D31" "0.0504"
$ws.Range("E31").Value = "  -1.32%  "
Set-TextValue "D32" "3.34"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "1.479.33"
$ws.Range("E33").Value = "  -0.92%  "
Set-TextValue "D34" "3.11"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -5.71%  "
Set-TextValue "D36" "0.935"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  -0.96%  "
Set-TextValue "D38" "0.0172"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("E39").Value = "  -6.16%  "
Set-TextValue "D40" "69.48"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("E41").Value = "  -5.78%  "
Set-TextValue "D42" "1.00"
$ws.Range("E42").Value = "  -0.03%  "
Set-TextValue "D43" "5.40"
$ws.Range("E43").Value = "  -7.46%  "
$ws.Range("D44").Value = "1.810.38"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("E47").Value = "  -2.45%  "
Set-TextValue "D48" "89.30"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("E50").Value = "  -2.71%  "
Set-TextValue "D51" "7.90"
$ws.Range("E51").Value = "  -3.44%  "
